$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.234.19"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "'2.883.07"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'525.83"
$ws.Range("E5").Value = "  -4.30%  "
$ws.Range("D6").Value = "'141.75"
$ws.Range("E6").Value = "  -6.80%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "'2.891.09"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").Value = "'3.386.89"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'60.286.29"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "'22.50"
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("D17").Value = "'2.880.39"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "'361.26"
$ws.Range("E21").Value = "  -7.59%  "
$ws.Range("D22").Value = "'6.52"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'63.24"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "'3.005.94"
$ws.Range("E25").Value = "  -4.09%  "
$ws.Range("D26").Value = "'0.448"
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("D27").Value = "'0.182"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'7.79"
$ws.Range("E29").Value = "  -7.62%  "
$ws.Range("D30").Value = "'0.0₃0852"
$ws.Range("E30").Value = "  -10.00%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "'19.43"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("D34").Value = "'148.61"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").Value = "'4.32"
$ws.Range("E35").Value = "  -7.27%  "
$ws.Range("E36").Value = "  -7.49%  "
$ws.Range("E37").Value = "  -7.41%  "
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").Value = "'37.74"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  -4.01%  "
$ws.Range("D41").Value = "'2.321.76"
$ws.Range("E41").Value = "  -4.95%  "
$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = "  -6.24%  "
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").Value = "'20.67"
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'5.05"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("D49").Value = "'10.35"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").Value = "'0.0928"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").Value = "'249.10"
$ws.Range("E51").Value = "  -5.02%  "
